# "Generate Report for Handoff"
#
# The b44cfec6-8057-4922-8a65-82d6d37b7f7b.md file has finished translation
# and is now ready to be handed off, for both the zh-cn and de-de locales.
# Update the per-locale detail sheets (Status/Priority/Latest Handoff
# Datetime) and roll the same status + timestamp up into the Overview
# sheet.

$wb = $excel.ActiveWorkbook

$newStatus  = "Ready for handoff"
$newPriority = "mt"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- zh-cn detail sheet: row 3 is b44cfec6-8057-4922-8a65-82d6d37b7f7b.md ---
$zhcn.Range("C3").Value2 = $newStatus
$zhcn.Range("E3").Value2 = $newPriority
$zhcn.Range("H3").Value2 = "2016-08-12 12:13:16"

# --- de-de detail sheet: row 3 is b44cfec6-8057-4922-8a65-82d6d37b7f7b.md ---
$dede.Range("C3").Value2 = $newStatus
$dede.Range("E3").Value2 = $newPriority
$dede.Range("H3").Value2 = "2016-08-12 12:13:23"

# --- Overview roll-up: row 3 is b44cfec6-8057-4922-8a65-82d6d37b7f7b.md ---
$overview.Range("E3").Value2 = $newStatus
$overview.Range("F3").Value2 = $newStatus
$overview.Range("G3").Value2 = "2016-08-12 12:13:23"

# Re-fit the Status columns now that "Ready for handoff" is wider than the
# previous "In Translation" text.
$zhcn.Columns.Item(3).AutoFit() | Out-Null
$dede.Columns.Item(3).AutoFit() | Out-Null
$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null
